$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2026-02-14 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-15 Sunday", 2)

# Update the division problems table, cell by cell (row, col are 1-based
# over the Word table model; only rows 1, 5, 9, 13, 17 contain data).
$t = $d.Tables.Item(1)

$values = @(
    @(1,  1, "91÷6=15, 1"),
    @(1,  2, "51÷2=25, 1"),
    @(1,  3, "42÷7=6, 0"),
    @(1,  4, "34÷8=4, 2"),
    @(1,  5, "32÷8=4, 0"),

    @(5,  1, "88÷9=9, 7"),
    @(5,  2, "61÷8=7, 5"),
    @(5,  3, "44÷7=6, 2"),
    @(5,  4, "49÷5=9, 4"),
    @(5,  5, "76÷7=10, 6"),

    @(9,  1, "53÷4=13, 1"),
    @(9,  2, "84÷5=16, 4"),
    @(9,  3, "56÷2=28, 0"),
    @(9,  4, "15÷7=2, 1"),
    @(9,  5, "19÷7=2, 5"),

    @(13, 1, "21÷9=2, 3"),
    @(13, 2, "24÷2=12, 0"),
    @(13, 3, "75÷7=10, 5"),
    @(13, 4, "66÷5=13, 1"),
    @(13, 5, "97÷8=12, 1"),

    @(17, 1, "72÷5=14, 2"),
    @(17, 2, "87÷3=29, 0"),
    @(17, 3, "49÷8=6, 1"),
    @(17, 4, "10÷5=2, 0"),
    @(17, 5, "25÷5=5, 0")
)

foreach ($v in $values) {
    $cell = $t.Cell($v[0], $v[1])
    $cell.Range.Text = $v[2]
}
